$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "dsf"
$ws.Range("K13").Value = "sdfdfsd"
$ws.Range("N6").Value = "sdf"

$ws.Range("N6").Select()
